{"js": "// Update the Montreal (CYUL) Tower emergency contact numbers:\n// \"514-633-3311 / 514-633-3312\"  ->  \"514-633-3312 (Primary) / 514-633-3311\"\nconst body = context.document.body;\nconst results = body.search(\"514-633-3311 / 514-633-3312\", { matchCase: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"514-633-3312 (Primary) / 514-633-3311\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the Montreal (CYUL) Tower emergency contact numbers:\n# \"514-633-3311 / 514-633-3312\"  ->  \"514-633-3312 (Primary) / 514-633-3311\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"514-633-3311 / 514-633-3312\"\n$found = $find.Execute()\n\nif ($found) {\n    $find.Parent.Text = \"514-633-3312 (Primary) / 514-633-3311\"\n}\n"}
